# Update countries & provincias Spain
#
# Updates the "Casos totales" (and related) figures for a handful of
# countries, then re-sorts the whole country table by "Casos totales"
# (column B) descending, which is how the sheet is always kept ordered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$countryCol = $ws.Range("A4:A218")

# New daily figures: Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes
$updates = @{
    "Austria"     = @(15558, 27, 13110, 1852, 124, 7, 596)
    "Polonia"     = @(13375, 270, 3762, 8949, 160, 13, 664)
    "Kazajistan"  = @(3800, 203, 940, 2835, 40, 0, 25)
    "Azerbaiyan"  = @(1894, 40, 1411, 458, 17, 0, 25)
    "Kuwait"      = @(4619, 242, 1703, 2883, 69, 3, 33)
    "Madagascar"  = @(135, 3, 97, 38, 1, 0, 0)
    "Zimbabue"    = @(34, 0, 5, 25, 0, 0, 4)
}

foreach ($name in $updates.Keys) {
    $found = $countryCol.Find($name)
    $r = $found.Row
    $vals = $updates[$name]
    $ws.Range("B" + $r).Value = $vals[0]
    $ws.Range("C" + $r).Value = $vals[1]
    $ws.Range("D" + $r).Value = $vals[2]
    $ws.Range("E" + $r).Value = $vals[3]
    $ws.Range("F" + $r).Value = $vals[4]
    $ws.Range("G" + $r).Value = $vals[5]
    $ws.Range("H" + $r).Value = $vals[6]
}

# Re-sort the whole data range (countries + values, no header) by
# "Casos totales" descending, as the sheet keeps it ranked.
$sortRange = $ws.Range("A4:H218")
$sortKey = $ws.Range("B4:B218")
$sortRange.Sort($sortKey, 2)
